# issue #5: property land done
#
# This script:
#  1) Fixes a batch of typo / stray-space / smart-quote data entry issues
#     that were scattered across every sheet of the property-disclosure
#     workbook (land, building, car, deposit, securities, insurance, debt).
#  2) Extends the "土地" (land) sheet with the new export-pipeline metadata
#     columns: property_category, category, date, legislator_name,
#     legislator_id, source_file, index - and renames its existing headers
#     to the new English field names used by the pipeline.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "土地" (Land) - sheet index 1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Existing headers get renamed to the pipeline's English field names.
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,3).Value = "area"
$ws.Cells.Item(1,4).Value = "share_portion"
$ws.Cells.Item(1,5).Value = "owner"
$ws.Cells.Item(1,6).Value = "register_date"
$ws.Cells.Item(1,7).Value = "register_reason"
$ws.Cells.Item(1,8).Value = "acquire_value"

# Fix the data-row typos: stray hyphen in the parcel number, stray space
# in the date.
$ws.Cells.Item(2,2).Value = "嘉義市長竹段00060002地號"
$ws.Cells.Item(2,6).Value = "96年02月09日"

# New trailing metadata columns (I:O) + header row.
$landHeaders = @("property_category","category","date","legislator_name","legislator_id","source_file","index")
for ($i = 0; $i -lt $landHeaders.Length; $i++) {
    $cell = $ws.Cells.Item(1, 9 + $i)
    $cell.Value = $landHeaders[$i]
    # Match the look of the other header cells (bold, centered, boxed).
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$ws.Cells.Item(2,9).Value  = "land"
$ws.Cells.Item(2,10).Value = "normal"
$ws.Cells.Item(2,11).Value = "'2012-04-24"
$ws.Cells.Item(2,12).Value = "李俊俋"
$ws.Cells.Item(2,13).Value = 1738
$ws.Cells.Item(2,14).Value = "tmp16861"
$ws.Cells.Item(2,15).Value = 15

# ---------------------------------------------------------------------
# Sheet "建物" (Building) - sheet index 2
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,2).Value = "嘉義市長竹段01946000建號"
$ws.Cells.Item(2,6).Value = "96年02月09日"

# ---------------------------------------------------------------------
# Sheet "汽車" (Car) - sheet index 3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(1,2).Value = "廠牌型號"
$ws.Cells.Item(2,2).Value = "MAZDA3"
$ws.Cells.Item(3,2).Value = "MAZDAMPV"
$ws.Cells.Item(3,6).Value = "maas貝賣"
$ws.Cells.Item(4,5).Value = "100年09月"

# ---------------------------------------------------------------------
# Sheet "存款" (Deposit) - sheet index 4
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2,2).Value  = "中華郵政股份有限公司台大郵局"
$ws.Cells.Item(5,2).Value  = "中華郵政股份有限公司嘉義中山路郵局"
$ws.Cells.Item(6,2).Value  = "玉山商業銀行"
$ws.Cells.Item(10,2).Value = "中華郵政股份有限公司嘉義中山路郵局"
$ws.Cells.Item(11,2).Value = "中華郵政股份有限公司嘉義中山路郵局"

# ---------------------------------------------------------------------
# Sheet "其他有價證券" (Other securities) - sheet index 5
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2,4).Value = "項"
$ws.Cells.Item(2,5).Value = "件|所有人"

# ---------------------------------------------------------------------
# Sheet "保險" (Insurance) - sheet index 6
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2,3).Value = "富邦人壽心得意利率變動型年金保險"

# ---------------------------------------------------------------------
# Sheet "債務" (Debt) - sheet index 7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2,6).Value = "96年02月12日"
